$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price data (fecha = 2021-09-10, serial 44449) is inserted just
# above the previous last two rows (formerly rows 66-67), which shift down
# to rows 70-71 unchanged.
$ws.Rows.Item(66).Resize(4).Insert()

function Set-PimientoRow {
    param($row, $calidad, $variedad, $volumen, $precio)

    $ws.Cells.Item($row, 1).Value = 12
    $ws.Cells.Item($row, 2).Value = "Mapocho Venta Directa de Santiago"
    $ws.Cells.Item($row, 3).Value = "Metropolitana"
    $ws.Cells.Item($row, 4).Value2 = 44449
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = 13
    $ws.Cells.Item($row, 6).Value = 100112002
    $ws.Cells.Item($row, 7).Value = "Pimiento"
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $calidad
    $ws.Cells.Item($row, 10).Value = $volumen
    $ws.Cells.Item($row, 11).Value = $precio
    $ws.Cells.Item($row, 12).Value = $precio
    $ws.Cells.Item($row, 13).Value = $precio
    $ws.Cells.Item($row, 14).Value = "$/caja 18 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 17).Value = 18
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

Set-PimientoRow 66 "Tercera" "Zafiro rojo" 25 34000
Set-PimientoRow 67 "Primera" "Zafiro verde" 15 38000
Set-PimientoRow 68 "Segunda" "Zafiro verde" 20 36000
Set-PimientoRow 69 "Tercera" "Zafiro verde" 25 34000

# Precio $/Kg (column P) values per row, computed from the source data
$ws.Cells.Item(66, 16).Value = 1889
$ws.Cells.Item(67, 16).Value = 2111
$ws.Cells.Item(68, 16).Value = 2000
$ws.Cells.Item(69, 16).Value = 1889
